# Generate Report for Handoff
#
# The localization status report moved from "In Translation" to
# "Ready for handoff": the Status cell on the zh-cn and de-de sheets is
# updated, along with the "Latest HO Xliff Generate Date" / "Latest
# Handoff Datetime" timestamps that get stamped when the handoff xliff is
# (re)generated. The Status column text got longer, so its column (and
# the two Overview columns that mirror each language's status) were
# widened to fit.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-04 19:03:39"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-04 19:03:43"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- Overview sheet: Latest HO Xliff Generate Date + mirrored status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 19:03:43"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
